# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment Schedule"
#   sheet, shifting the old N/O/P ("Late"/heading/"Outstanding") data right
#   to O/P/Q.
# - Make "Repayment Schedule" the active sheet/tab, with M12 selected.
# - "Modify Transaction" is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$wsModifyTransaction = $wb.Worksheets.Item("Modify Transaction")
$wsRepaymentSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column before column N (pushes existing N:P -> O:Q).
$wsRepaymentSchedule.Columns("N:N").Insert()

# Move the view: Repayment Schedule becomes the active/selected sheet,
# Modify Transaction is no longer tabSelected.
$wsRepaymentSchedule.Activate()
$wsRepaymentSchedule.Range("M12").Select()

Write-Output "done"
